$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextCell 2 4 '37.196.67'
Set-TextCell 2 5 '  +0.34%  '
Set-TextCell 3 4 '2.058.56'
Set-TextCell 3 5 '  +0.17%  '
Set-TextCell 4 5 '  -0.14%  '
Set-TextCell 5 4 '248.75'
Set-TextCell 5 5 '  -1.47%  '
Set-TextCell 6 5 '  -0.59%  '
Set-TextCell 7 5 '  -0.05%  '
Set-TextCell 8 4 '57.32'
Set-TextCell 8 5 '  -2.34%  '
Set-TextCell 9 5 '  +0.16%  '
Set-TextCell 10 4 '0.0787'
Set-TextCell 10 5 '  -1.34%  '
Set-TextCell 11 5 '  +0.29%  '
Set-TextCell 12 4 '16.33'
Set-TextCell 12 5 '  -0.46%  '
Set-TextCell 13 4 '0.922'
Set-TextCell 13 5 '  +14.37%  '
Set-TextCell 14 4 '2.355.25'
Set-TextCell 14 5 '  -0.08%  '
Set-TextCell 15 4 '5.81'
Set-TextCell 15 5 '  +3.74%  '
Set-TextCell 16 4 '2.059.97'
Set-TextCell 16 5 '  -0.31%  '
Set-TextCell 17 4 '18.83'
Set-TextCell 17 5 '  +13.68%  '
Set-TextCell 18 4 '37.196.18'
Set-TextCell 18 5 '  +0.51%  '
Set-TextCell 19 4 '75.00'
Set-TextCell 19 5 '  -0.96%  '
Set-TextCell 20 4 '0.0₃0903'
Set-TextCell 20 5 '  -1.68%  '
Set-TextCell 21 4 '5.51'
Set-TextCell 21 5 '  +0.76%  '
Set-TextCell 22 4 '238.12'
Set-TextCell 22 5 '  -0.03%  '
Set-TextCell 23 5 '  +0.00%  '
Set-TextCell 24 5 '  +4.44%  '
Set-TextCell 25 5 '  +4.45%  '
Set-TextCell 26 4 '2.20'
Set-TextCell 26 5 '  -4.34%  '
Set-TextCell 27 4 '170.37'
Set-TextCell 27 5 '  +0.69%  '
Set-TextCell 28 4 '20.30'
Set-TextCell 28 5 '  +0.13%  '
Set-TextCell 29 5 '  -0.84%  '
Set-TextCell 30 4 '5.19'
Set-TextCell 30 5 '  +9.14%  '
Set-TextCell 31 4 '1.17'
Set-TextCell 31 5 '  +2.19%  '
Set-TextCell 32 5 '  +0.68%  '
Set-TextCell 33 4 '4.66'
Set-TextCell 33 5 '  +3.87%  '
Set-TextCell 34 4 '0.0886'
Set-TextCell 34 5 '  -0.32%  '
Set-TextCell 35 5 '  -0.05%  '
Set-TextCell 36 4 '2.29'
Set-TextCell 36 5 '  +0.43%  '
Set-TextCell 37 4 '1.77'
Set-TextCell 37 5 '  +0.98%  '
Set-TextCell 38 5 '  -1.50%  '
Set-TextCell 39 4 '5.29'
Set-TextCell 39 5 '  +17.28%  '
Set-TextCell 40 5 '  +7.52%  '
Set-TextCell 41 4 '0.100'
Set-TextCell 41 5 '  -13.16%  '
Set-TextCell 42 4 '17.76'
Set-TextCell 42 5 '  -0.16%  '
Set-TextCell 44 5 '  +0.43%  '
Set-TextCell 45 4 '96.89'
Set-TextCell 45 5 '  -0.77%  '
Set-TextCell 46 5 '  -1.33%  '
Set-TextCell 47 4 '1.278.59'
Set-TextCell 47 5 '  -1.08%  '
Set-TextCell 48 4 '2.86'
Set-TextCell 48 5 '  -1.63%  '
Set-TextCell 49 5 '  -0.47%  '
Set-TextCell 50 4 '2.243.73'
Set-TextCell 50 5 '  -0.16%  '
Set-TextCell 51 4 '44.46'
Set-TextCell 51 5 '  +0.96%  '
